# Fruta / hortaliza, semanal
# Update the weekly price records: dates and volume/price figures were
# re-shuffled across rows 2,3,5,6,7,8,9,10,11 (row 4 is unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44533
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 16000
$ws.Range("O2").Value = 17000
$ws.Range("P2").Value = 16500
$ws.Range("S2").Value = 825

# Row 3
$ws.Range("D3").Value = 44320
$ws.Range("M3").Value = 80
$ws.Range("N3").Value = 16000
$ws.Range("O3").Value = 17000
$ws.Range("P3").Value = 16500
$ws.Range("S3").Value = 825

# Row 5
$ws.Range("D5").Value = 44357
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 14500
$ws.Range("S5").Value = 725

# Row 6
$ws.Range("D6").Value = 45092
$ws.Range("M6").Value = 150
$ws.Range("N6").Value = 24000
$ws.Range("O6").Value = 25000
$ws.Range("P6").Value = 24333
$ws.Range("S6").Value = 1217

# Row 7
$ws.Range("D7").Value = 44893
$ws.Range("M7").Value = 80
$ws.Range("N7").Value = 21000
$ws.Range("O7").Value = 22000
$ws.Range("P7").Value = 21625
$ws.Range("S7").Value = 1081

# Row 8
$ws.Range("D8").Value = 44708
$ws.Range("M8").Value = 80
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 21000
$ws.Range("P8").Value = 20500
$ws.Range("S8").Value = 1025

# Row 9
$ws.Range("D9").Value = 44761
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 21000
$ws.Range("P9").Value = 20500
$ws.Range("S9").Value = 1025

# Row 10
$ws.Range("D10").Value = 44890
$ws.Range("M10").Value = 80
$ws.Range("N10").Value = 20000
$ws.Range("O10").Value = 23000
$ws.Range("P10").Value = 22250
$ws.Range("S10").Value = 1112

# Row 11
$ws.Range("D11").Value = 44792
$ws.Range("M11").Value = 100
$ws.Range("N11").Value = 21000
$ws.Range("O11").Value = 22000
$ws.Range("P11").Value = 21500
$ws.Range("S11").Value = 1075
